$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2: Target cluster changes from MuSCs to ECs, and
# several downstream metrics (M..T) are recomputed with the new TPM values. ---
$ws.Range("D2").Value = "ECs"

$ws.Range("M2").Value = 0.044174
$ws.Range("N2").Value = 0.132522
$ws.Range("O2").Value = 0.2474193313505733
$ws.Range("P2").Value = 0.2474193313505733
$ws.Range("Q2").Value = 0.07164265952133332
$ws.Range("R2").Value = 0.644783935692
$ws.Range("S2").Value = 0.2474193313505733
$ws.Range("T2").Value = 0.2474193313505733

# --- Add new row 3: the (FAPs, Rspo2, Lgr6, MuSCs) LR-pair row that used to
# be represented by row 2 before the MuSCs/ECs split, now recomputed too. ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "MuSCs"

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.621828666666667
$ws.Range("H3").Value = 4.865486
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.134365
$ws.Range("N3").Value = 0.403095
$ws.Range("O3").Value = 0.7525806686494267
$ws.Range("P3").Value = 0.7525806686494266
$ws.Range("Q3").Value = 0.2179170087966666
$ws.Range("R3").Value = 1.96125307917
$ws.Range("S3").Value = 0.7525806686494267
$ws.Range("T3").Value = 0.7525806686494266
